$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the placeholder image-link text from rows that no longer need it
$ws.Range("F1").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("F11").ClearContents()

# Update remaining rows with their real logo links
$ws.Range("F2").Value = "https://cdn.prod.website-files.com/657c80cc477a31a499849a7f/667f1639c800c15e37c46c78_barcode-generator-logo.webp"
$ws.Range("F4").Value = "https://upload.wikimedia.org/wikipedia/fr/thumb/2/2e/R%C3%A9gion_Hauts-de-France_logo_2016.svg/2048px-R%C3%A9gion_Hauts-de-France_logo_2016.svg.png"

$ws.Range("F5").Value = "https://s3-eu-west-1.amazonaws.com/tpd/logos/5447a9cf00006400057b13a3/0x0.png"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://s3-eu-west-1.amazonaws.com/tpd/logos/5447a9cf00006400057b13a3/0x0.png") | Out-Null

$ws.Range("F8").Value = "https://png.pngtree.com/png-vector/20230423/ourmid/pngtree-appointment-line-icon-vector-png-image_6720015.png"
$ws.Range("F9").Value = "https://www.verifiance-fnci.fr/Souscription/assets/ctx/2374d0af/resources/aem-import/css/images/logo-verifiance.png"
$ws.Range("F13").Value = "https://leparcduluc.fr/wp-content/uploads/2022/08/1200px-Kiabi_logo.svg.png"
$ws.Range("F14").Value = "https://creation-entreprise.info/wp-content/uploads/2024/05/smart-rh.jpg"

# Move the active selection to F14, matching the saved cursor position
$ws.Range("F14").Select()
